$wb = $excel.ActiveWorkbook

$sheetNames = @(
    "Regular_Section_A",
    "Regular_Section_B",
    "PreMid_Section_A",
    "PreMid_Section_B",
    "PostMid_Section_A",
    "PostMid_Section_B"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("B2").Value = "MINOR: Generative Ai [C102]"
    $ws.Range("C8").Value = "CS303 (Lab) [L107]"
    $ws.Range("C9").Value = "CS303 (Lab) [L107]"
    $ws.Range("B10").Value = "MINOR: VLSI [C102]"
}
